$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Daryl Dixon) "Orig Pwd" cell changes from "Labo0749" to "P@ssw0rd2".
$ws.Range("C3").Value = "P@ssw0rd2"

# Excel auto-hyperlinked the "@"-containing cells in row 3, the same way the
# existing C2 cell is already hyperlinked to "mailto:P@ssw0rd1" (matching its
# own text). Add the same mailto hyperlinks for A3 (the agent's email
# address) and C3 (the new password).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:DarylD@TestIncidentQueue.onmicrosoft.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:P@ssw0rd2") | Out-Null

# Re-apply the same built-in Hyperlink cell style already used on C2 so the
# two newly-linked cells look consistent with it.
$ws.Range("A3").Style = $ws.Range("C2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style

# The active selection moved to D3.
$ws.Range("D3").Select() | Out-Null
